$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.904.20'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.634.66'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5071'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2578'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E9').Value = '  -0.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.62'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.83%  '
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.259'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.642.46'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.07%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5512'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0₅7697'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.89'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.918.05'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.437'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '194.92'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.912'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.049'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.003'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.902'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '142.29'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1241'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.810'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.60'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.241'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.04890'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.251'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.190'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.541'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.19%  '
$ws.Range('E34').Value = '  +0.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9053'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.16%  '
$ws.Range('E36').Value = '  -1.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5501'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.124.57'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.84%  '
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('E40').Value = '  -0.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.571'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8021'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '97.59'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.92%  '
$ws.Range('E44').Value = '  -4.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.773.47'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4453'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.83%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.79'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9951'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.72%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05150'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.548'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.004'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.11%  '
